$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 784, shifting the existing rows 784..825 down to 785..826
$ws.Rows.Item(784).Insert()

# The new row represents another reading for 2026/02/12 (Thursday) at hour 13
# Make sure the date/day-of-week are written as plain text (matching the rest of
# column A/B which store the date strings as text, not as Excel date serials).
$ws.Range("A784:B784").NumberFormat = "@"
$ws.Range("A784").Value = "2026/02/12"
$ws.Range("B784").Value = "木"
$ws.Range("C784").Value = 13
$ws.Range("D784").Value = 27

# Drop the temporary text format override so the new cells fall back to the
# workbook's default (unstyled) look, just like every other data row.
$ws.Range("A784:B784").Style = "Normal"

Write-Output "done"
